$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3826.2666
$ws.Range("I74").Value = 3616.1667
$ws.Range("K74").Value = 3616.1667
$ws.Range("M74").Value = -2680.1667
$ws.Range("H77").Value = 3826.2666
$ws.Range("I77").Value = 3616.1667
$ws.Range("K77").Value = 18080.8335
$ws.Range("M77").Value = -13400.8335
$ws.Range("H92").Value = 541.8
$ws.Range("J92").Value = 999
$ws.Range("L92").Value = 999
$ws.Range("N92").Value = -3495
$ws.Range("H99").Value = 572.7273
$ws.Range("I99").Value = 575.125
$ws.Range("K99").Value = 1725.375
$ws.Range("M99").Value = -227.375
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H108").Value = 56000
$ws.Range("I108").Value = 48000
$ws.Range("J108").Value = 60000
$ws.Range("K108").Value = 48000
$ws.Range("L108").Value = 60000
$ws.Range("M108").Value = -44160
$ws.Range("N108").Value = -67680
$ws.Range("H113").Value = 5099.5713
$ws.Range("I113").Value = 4235
$ws.Range("J113").Value = 5748
$ws.Range("K113").Value = 4235
$ws.Range("L113").Value = 5748
$ws.Range("M113").Value = -981
$ws.Range("N113").Value = -12256

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 186
$ws.Range("I5").Value = 158.5
$ws.Range("K5").Value = 158.5
$ws.Range("M5").Value = -46.5
$ws.Range("H45").Value = 31252292
$ws.Range("J45").Value = 3338
$ws.Range("L45").Value = 3338
$ws.Range("N45").Value = -4092
$ws.Range("H61").Value = 37504480
$ws.Range("J61").Value = 83340000
$ws.Range("L61").Value = 83340000
$ws.Range("N61").Value = -83340424
$ws.Range("H102").Value = 5195.8184
$ws.Range("I102").Value = 5195.8184
$ws.Range("K102").Value = 5195.8184
$ws.Range("M102").Value = -3573.8184
$ws.Range("H110").Value = 1533.1765
$ws.Range("I110").Value = 1262.3636
$ws.Range("K110").Value = 1262.3636
$ws.Range("M110").Value = 782.6364000000001
$ws.Range("H119").Value = 70174
$ws.Range("J119").Value = 70174
$ws.Range("L119").Value = 70174
$ws.Range("N119").Value = -79850
$ws.Range("H136").Value = 37504480
$ws.Range("J136").Value = 83340000
$ws.Range("L136").Value = 250020000
$ws.Range("N136").Value = -250025100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 186
$ws.Range("I4").Value = 158.5
$ws.Range("K4").Value = 158.5
$ws.Range("M4").Value = -43.5
$ws.Range("H20").Value = 5299.0605
$ws.Range("I20").Value = 4901.3
$ws.Range("K20").Value = 4901.3
$ws.Range("M20").Value = -4654.3
$ws.Range("H22").Value = 2538.7273
$ws.Range("I22").Value = 1560.8572
$ws.Range("J22").Value = 4250
$ws.Range("K22").Value = 1560.8572
$ws.Range("L22").Value = 4250
$ws.Range("M22").Value = -1387.8572
$ws.Range("N22").Value = -4596
$ws.Range("H86").Value = 2939.037
$ws.Range("I86").Value = 3078.3
$ws.Range("J86").Value = 2541.1428
$ws.Range("K86").Value = 3078.3
$ws.Range("L86").Value = 2541.1428
$ws.Range("M86").Value = -1955.3
$ws.Range("N86").Value = -4787.1428
$ws.Range("H89").Value = 2939.037
$ws.Range("I89").Value = 3078.3
$ws.Range("J89").Value = 2541.1428
$ws.Range("K89").Value = 15391.5
$ws.Range("L89").Value = 12705.714
$ws.Range("M89").Value = -9775.5
$ws.Range("N89").Value = -23937.714
$ws.Range("H105").Value = 2134.8857
$ws.Range("J105").Value = 2785.9473
$ws.Range("L105").Value = 2785.9473
$ws.Range("N105").Value = -6279.9473
$ws.Range("H107").Value = 8001.8
$ws.Range("I107").Value = 8001.8
$ws.Range("K107").Value = 8001.8
$ws.Range("M107").Value = -6081.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 943154.1
$ws.Range("I31").Value = 14984.728
$ws.Range("J31").Value = 1672430.1
$ws.Range("K31").Value = 14984.728
$ws.Range("L31").Value = 1672430.1
$ws.Range("M31").Value = -14689.728
$ws.Range("N31").Value = -1673020.1
$ws.Range("H34").Value = 943154.1
$ws.Range("I34").Value = 14984.728
$ws.Range("J34").Value = 1672430.1
$ws.Range("K34").Value = 14984.728
$ws.Range("L34").Value = 1672430.1
$ws.Range("M34").Value = -14782.728
$ws.Range("N34").Value = -1672834.1
$ws.Range("H132").Value = 7430.0586
$ws.Range("I132").Value = 2789.077
$ws.Range("K132").Value = 8367.231
$ws.Range("M132").Value = -5837.231

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2115.6667
$ws.Range("J5").Value = 2500
$ws.Range("L5").Value = 7500
$ws.Range("N5").Value = -7724
$ws.Range("H34").Value = 601
$ws.Range("I34").Value = 268
$ws.Range("J34").Value = 934
$ws.Range("K34").Value = 804
$ws.Range("L34").Value = 2802
$ws.Range("M34").Value = -720
$ws.Range("N34").Value = -2970
$ws.Range("H59").Value = 1780
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1780
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5340
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -6420
$ws.Range("H107").Value = 984.5833
$ws.Range("J107").Value = 1131.5
$ws.Range("L107").Value = 3394.5
$ws.Range("N107").Value = -7234.5
$ws.Range("H135").Value = 2115.6667
$ws.Range("J135").Value = 2500
$ws.Range("L135").Value = 22500
$ws.Range("N135").Value = -27570
$ws.Range("H139").Value = 6714.222
$ws.Range("I139").Value = 8857.5
$ws.Range("J139").Value = 4999.6
$ws.Range("K139").Value = 26572.5
$ws.Range("L139").Value = 14998.8
$ws.Range("M139").Value = -21432.5
$ws.Range("N139").Value = -25278.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 128341.2
$ws.Range("J108").Value = 128341.2
$ws.Range("L108").Value = 128341.2
$ws.Range("N108").Value = -136021.2
$ws.Range("H113").Value = 3479.5789
$ws.Range("I113").Value = 2767.2222
$ws.Range("J113").Value = 4120.7
$ws.Range("K113").Value = 2767.2222
$ws.Range("L113").Value = 4120.7
$ws.Range("M113").Value = -597.2222000000002
$ws.Range("N113").Value = -8460.700000000001
$ws.Range("H126").Value = 3751.6
$ws.Range("I126").Value = 3068.8125
$ws.Range("J126").Value = 4965.4443
$ws.Range("K126").Value = 9206.4375
$ws.Range("L126").Value = 14896.3329
$ws.Range("M126").Value = -6736.4375
$ws.Range("N126").Value = -19836.3329
$ws.Range("H132").Value = 28573680
$ws.Range("I132").Value = 35716410
$ws.Range("K132").Value = 107149230
$ws.Range("M132").Value = -107146700

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2439.4
$ws.Range("I22").Value = 2291.7693
$ws.Range("J22").Value = 2713.5715
$ws.Range("K22").Value = 2291.7693
$ws.Range("L22").Value = 2713.5715
$ws.Range("M22").Value = -1996.7693
$ws.Range("N22").Value = -3303.5715
$ws.Range("H27").Value = 2439.4
$ws.Range("I27").Value = 2291.7693
$ws.Range("J27").Value = 2713.5715
$ws.Range("K27").Value = 2291.7693
$ws.Range("L27").Value = 2713.5715
$ws.Range("M27").Value = -2184.7693
$ws.Range("N27").Value = -2927.5715
$ws.Range("H121").Value = 113765
$ws.Range("J121").Value = 113765
$ws.Range("L121").Value = 113765
$ws.Range("N121").Value = -117259
$ws.Range("H132").Value = 582652.5600000001
$ws.Range("I132").Value = 59466.332
$ws.Range("K132").Value = 178398.996
$ws.Range("M132").Value = -175868.996

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3110.5
$ws.Range("I96").Value = 2294.75
$ws.Range("J96").Value = 3926.25
$ws.Range("K96").Value = 2294.75
$ws.Range("L96").Value = 3926.25
$ws.Range("M96").Value = -921.75
$ws.Range("N96").Value = -6672.25
$ws.Range("H100").Value = 546.5833
$ws.Range("I100").Value = 438.6111
$ws.Range("J100").Value = 870.5
$ws.Range("K100").Value = 877.2222
$ws.Range("L100").Value = 1741
$ws.Range("M100").Value = -336.2222
$ws.Range("N100").Value = -2823
$ws.Range("H107").Value = 26317138
$ws.Range("I107").Value = 33334522
$ws.Range("J107").Value = 1948
$ws.Range("K107").Value = 100003566
$ws.Range("L107").Value = 5844
$ws.Range("M107").Value = -100001646
$ws.Range("N107").Value = -9684
$ws.Range("H113").Value = 691.26086
$ws.Range("J113").Value = 930.75
$ws.Range("L113").Value = 2792.25
$ws.Range("N113").Value = -7132.25
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H126").Value = 2099.75
$ws.Range("I126").Value = 2133
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6399
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3929
$ws.Range("N126").Value = -10940
